$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(25).Delete()
[void]$ws.Rows.Item(25).Select()
